# Regenerate merged AHB files
# - rename header labels: "_old" -> "_FV2310", "_new" -> "_FV2404"
# - turn the data range into an Excel Table (ListObject) with autofilter
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row labels (row 1, columns A:J and L:U carry the
#    "_old"/"_new" suffixed headers; K holds the literal "diff" header).
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '_old$', '_FV2310'
        $newVal = $newVal -replace '_new$', '_FV2404'
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

# 2. Convert the used range A1:U57 into a native Excel Table with an
#    autofilter, matching the tableParts/table1.xml addition.
$tableRange = $ws.Range("A1:U57")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# 3. Freeze the header row (pane split after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
